# A new daily price record (2022-03-03, serial 44623) is inserted at row 36
# of the sheet, which pushes all existing records from row 36 down through
# row 151 to rows 37 through 152 (dimension grows from A1:R151 to A1:R152).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 36; Excel shifts row 36..151 down to 37..152
# and extends the used range automatically.
$ws.Rows.Item(36).Insert()

# Fill in the new record's data in the now-empty row 36.
$ws.Range("A36").Value = 3
$ws.Range("B36").Value = "Femacal de La Calera"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 44623
$ws.Range("E36").Value = 5
$ws.Range("F36").Value = 100112052
$ws.Range("G36").Value = "Albahaca"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 105
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 5500
$ws.Range("M36").Value = 5262
$ws.Range("N36").Value = "`$/docena de matas"
$ws.Range("O36").Value = "Provincia de Quillota"
$ws.Range("P36").Value = 877
$ws.Range("Q36").Value = 6
$ws.Range("R36").Value = "Hortaliza"
